$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mark ABS() (row 2) and ASC() (row 3) as supported (Token, Parse, Eval) with "X",
# matching the existing pattern used for other supported functions (e.g. row 5 AUTO).
$ws.Range("B2").Value = "X"
$ws.Range("C2").Value = "X"
$ws.Range("D2").Value = "X"

$ws.Range("B3").Value = "X"
$ws.Range("C3").Value = "X"
$ws.Range("D3").Value = "X"

# LEN() (row 76) gains an "Eval" mark too (D column).
$ws.Range("D76").Value = "X"

# Move the active selection to D3, matching the updated selection in the sheet view.
$ws.Range("D3").Select()
